$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "send http post request" row with a new URL (custom headers feature).
$ws.Range("B6").Value = "http://requestb.in/1jv3zq61"

# Add a new "headers" column: label in row 5 (next to the other field labels)
# and the actual header value in row 6 (next to the other field values).
$ws.Range("F5").Value = "headers"
$ws.Range("F6").Value = "User-Agent=testx;Something='else entirely'"

# Restore the view so the top-left visible cell is A1 and the active cell is F5.
$ws.Range("A1").Select()
$ws.Range("F5:F6").Select()
